$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the pooled standard error formulas: the denominator should always
# reference the fixed row 5 (C$5) instead of incorrectly shifted cells.
$ws.Range("F2").Formula = "=SQRT(((C2-1)*E2^2+(C`$5-1)*E`$5^2)/(C2+C`$5-2))"
$ws.Range("F3").Formula = "=SQRT(((C3-1)*E3^2+(C`$5-1)*E`$5^2)/(C3+C`$5-2))"
$ws.Range("F4").Formula = "=SQRT(((C4-1)*E4^2+(C`$5-1)*E`$5^2)/(C4+C`$5-2))"

# Update the active cell selection to F5
$ws.Range("F5").Select()
